# Locate the paragraph containing the sentence that needs fixing:
#   "... Shouldn't be to large an obstacle."
# and turn it into:
#   "... Shouldn't be too large an obstacle."
# while also dropping the now-stale grammar-check markers
# (<w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>)
# that used to flag the word "to" as a likely grammar mistake.

$d = $word.ActiveDocument

$apostrophe = [char]0x2019

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t -like ("*Shouldn" + $apostrophe + "t be to large an obstacle*")) {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # Rebuild the paragraph's OOXML exactly as before, except the
    # "to" run loses its grammar-error wrapper and becomes "too".
    $ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' " +
          "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

    $paraXml =
        "<w:p $ns w14:paraId=`"49559946`" w14:textId=`"7CAD7629`" w:rsidR=`"000B587E`" w:rsidRDefault=`"006727E3`" w:rsidP=`"000B587E`">" +
        "<w:r><w:t xml:space=`"preserve`">Understanding forms.py and what Jake made in the </w:t></w:r>" +
        "<w:r w:rsidR=`"00601036`"><w:t>high-fidelity</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`"> prototype.</w:t></w:r>" +
        "<w:r w:rsidR=`"009C7096`"><w:t xml:space=`"preserve`"> Shouldn" + $apostrophe + "t be </w:t></w:r>" +
        "<w:r w:rsidR=`"009C7096`"><w:t>too</w:t></w:r>" +
        "<w:r w:rsidR=`"009C7096`"><w:t xml:space=`"preserve`"> large an obstacle.</w:t></w:r>" +
        "</w:p>"

    $r.InsertXML($paraXml)
} else {
    Write-Output "target paragraph not found"
}
